# Mejorando la version de la app
$wb = $excel.ActiveWorkbook

# --- Sheet "Productos": update rows 2-3, remove row 4 ---
$wsProductos = $wb.Worksheets.Item("Productos")

$wsProductos.Range("A2").Value = "Diadema"
$wsProductos.Range("B2").Value = "POLY 45EU36"

$wsProductos.Range("A3").Value = "MouseAlambrico"
$wsProductos.Range("B3").Value = "Mouse HP Alambrico"

$wsProductos.Rows.Item(4).Delete()

# --- Sheet "Proveedores": append rows 5-7 ---
$wsProveedores = $wb.Worksheets.Item("Proveedores")

$wsProveedores.Range("A5").Value = "Mauricio"
$wsProveedores.Range("B5").Value = "maoma.gaviria@hotmail.com"

$wsProveedores.Range("A6").Value = "Admin"
$wsProveedores.Range("B6").Value = "damiangaviria8@gmail.com"

$wsProveedores.Range("A7").Value = "Admin2"
$wsProveedores.Range("B7").Value = "damian.gaviria@est.iudigital.edu.co"
